# Update gh-pages to output generated at 456a3b4
# Applies updated "F" column (collected / sold count) values on the
# "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 2822
$ws1.Range("F11").Value = 487
$ws1.Range("F13").Value = 424
$ws1.Range("F19").Value = 871
$ws1.Range("F20").Value = 872
$ws1.Range("F26").Value = 8
$ws1.Range("F29").Value = 41
$ws1.Range("F31").Value = 1600
$ws1.Range("F32").Value = 377
$ws1.Range("F34").Value = 1529
$ws1.Range("F36").Value = 2369
$ws1.Range("F39").Value = 622
$ws1.Range("F44").Value = 1502
$ws1.Range("F45").Value = 221
$ws1.Range("F48").Value = 56

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 2822
$ws4.Range("F8").Value  = 487
$ws4.Range("F10").Value = 424
$ws4.Range("F15").Value = 872
$ws4.Range("F16").Value = 872
$ws4.Range("F27").Value = 41
$ws4.Range("F29").Value = 1600
$ws4.Range("F30").Value = 377
$ws4.Range("F33").Value = 2369
$ws4.Range("F39").Value = 622
$ws4.Range("F44").Value = 1502
$ws4.Range("F46").Value = 221
$ws4.Range("F48").Value = 56
